$wb = $excel.ActiveWorkbook

# --- Update selection on "Não é realizar a rota" sheet (5th sheet) ---
# Target diff changes its <selection> from sqref="A1:XFD1048576" to sqref="A1:E16"
$wsRota = $wb.Worksheets.Item("Não é realizar a rota")
$wsRota.Range("A1:E16").Select() | Out-Null

# --- Add the new worksheet "Planilha1" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Planilha1"

# Header row (reuses existing shared strings "#","Endereço","Abertura","Fechamento","Descarga")
$ws.Range("A1").Value = "#"
$ws.Range("B1").Value = "Endereço"
$ws.Range("C1").Value = "Abertura"
$ws.Range("D1").Value = "Fechamento"
$ws.Range("E1").Value = "Descarga"
$ws.Range("A1:E1").Interior.Color = 13553360

# Row 2 (index 0)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Metrô Consolação - Avenida Paulista - Cerqueira César, São Paulo - SP, Brasil"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("C2:D2").NumberFormat = "h:mm"
$ws.Range("C2").Interior.Color = 15132391

# Row 3 (index 1) - shaded row
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Manaus, AM, Brasil"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("A3:E3").Interior.Color = 15132391
$ws.Range("C3:D3").NumberFormat = "h:mm"

# Row 4 (index 2)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Minas Gerais, Brasil"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("C4:D4").NumberFormat = "h:mm"
$ws.Range("C4").Interior.Color = 15132391

# Column B width (auto-fit style width for the address column)
$ws.Columns.Item(2).ColumnWidth = 67

# Narrow / metric page margins matching this workbook's other sheets' locale defaults
$ws.PageSetup.LeftMargin = 0.511811024 * 72
$ws.PageSetup.RightMargin = 0.511811024 * 72
$ws.PageSetup.TopMargin = 0.78740157499999996 * 72
$ws.PageSetup.BottomMargin = 0.78740157499999996 * 72
$ws.PageSetup.HeaderMargin = 0.31496062000000002 * 72
$ws.PageSetup.FooterMargin = 0.31496062000000002 * 72

# Final selection / active cell on the new sheet
$ws.Range("B10").Select() | Out-Null
